$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift the footer (signature block) rows down by one to make room
#        for the new data row appended to the table. ---
$ws.Rows.Item(129).Insert()

# --- 2. The table used to end at row 124 with a special "closing border"
#        style; extend the table by one row (new row 125 gets that closing
#        style) and make the old row 124 a normal interior row. ---
$ws.Range("B124:J124").Copy()
$ws.Range("B125:J125").PasteSpecial(-4122)

$ws.Range("B123:J123").Copy()
$ws.Range("B124:J124").PasteSpecial(-4122)

# --- 3. Header updates ---
$ws.Range("E11").Value = 3265556
$ws.Range("F13").Value = 110

# --- 4. Rewrite the data table in ascending period order (previously
#        descending) and append the new period 2508 as the last row. ---
$ws.Range("E16").Value = "1607"
$ws.Range("F16").Value = 24640
$ws.Range("E17").Value = "1608"
$ws.Range("F17").Value = 24640
$ws.Range("E18").Value = "1609"
$ws.Range("F18").Value = 24640
$ws.Range("E19").Value = "1610"
$ws.Range("F19").Value = 24640
$ws.Range("E20").Value = "1611"
$ws.Range("F20").Value = 24640
$ws.Range("E21").Value = "1612"
$ws.Range("F21").Value = 24640
$ws.Range("E22").Value = "1701"
$ws.Range("F22").Value = 24640
$ws.Range("E23").Value = "1702"
$ws.Range("F23").Value = 24640
$ws.Range("E24").Value = "1703"
$ws.Range("F24").Value = 24640
$ws.Range("E25").Value = "1704"
$ws.Range("F25").Value = 24640
$ws.Range("E26").Value = "1705"
$ws.Range("F26").Value = 24640
$ws.Range("E27").Value = "1706"
$ws.Range("F27").Value = 24640
$ws.Range("E28").Value = "1707"
$ws.Range("F28").Value = 24640
$ws.Range("E29").Value = "1708"
$ws.Range("F29").Value = 24640
$ws.Range("E30").Value = "1709"
$ws.Range("F30").Value = 24640
$ws.Range("E31").Value = "1710"
$ws.Range("F31").Value = 24640
$ws.Range("E32").Value = "1711"
$ws.Range("F32").Value = 24640
$ws.Range("E33").Value = "1712"
$ws.Range("F33").Value = 24640
$ws.Range("E34").Value = "1801"
$ws.Range("F34").Value = 24640
$ws.Range("E35").Value = "1802"
$ws.Range("F35").Value = 24640
$ws.Range("E36").Value = "1803"
$ws.Range("F36").Value = 24640
$ws.Range("E37").Value = "1804"
$ws.Range("F37").Value = 24640
$ws.Range("E38").Value = "1805"
$ws.Range("F38").Value = 24640
$ws.Range("E39").Value = "1806"
$ws.Range("F39").Value = 24640
$ws.Range("E40").Value = "1807"
$ws.Range("F40").Value = 24640
$ws.Range("E41").Value = "1808"
$ws.Range("F41").Value = 24640
$ws.Range("E42").Value = "1809"
$ws.Range("F42").Value = 31249
$ws.Range("E43").Value = "1810"
$ws.Range("F43").Value = 31249
$ws.Range("E44").Value = "1811"
$ws.Range("F44").Value = 31249
$ws.Range("E45").Value = "1812"
$ws.Range("F45").Value = 31249
$ws.Range("E46").Value = "1901"
$ws.Range("F46").Value = 31249
$ws.Range("E47").Value = "1902"
$ws.Range("F47").Value = 31249
$ws.Range("E48").Value = "1903"
$ws.Range("F48").Value = 31249
$ws.Range("E49").Value = "1904"
$ws.Range("F49").Value = 31249
$ws.Range("E50").Value = "1905"
$ws.Range("F50").Value = 31249
$ws.Range("E51").Value = "1906"
$ws.Range("F51").Value = 31249
$ws.Range("E52").Value = "1907"
$ws.Range("F52").Value = 31249
$ws.Range("E53").Value = "1908"
$ws.Range("F53").Value = 31249
$ws.Range("E54").Value = "1909"
$ws.Range("F54").Value = 31249
$ws.Range("E55").Value = "1910"
$ws.Range("F55").Value = 31249
$ws.Range("E56").Value = "1911"
$ws.Range("F56").Value = 31249
$ws.Range("E57").Value = "1912"
$ws.Range("F57").Value = 31249
$ws.Range("E58").Value = "2001"
$ws.Range("F58").Value = 31249
$ws.Range("E59").Value = "2002"
$ws.Range("F59").Value = 31249
$ws.Range("E60").Value = "2003"
$ws.Range("F60").Value = 31249
$ws.Range("E61").Value = "2004"
$ws.Range("F61").Value = 31249
$ws.Range("E62").Value = "2005"
$ws.Range("F62").Value = 31249
$ws.Range("E63").Value = "2006"
$ws.Range("F63").Value = 31249
$ws.Range("E64").Value = "2007"
$ws.Range("F64").Value = 31249
$ws.Range("E65").Value = "2008"
$ws.Range("F65").Value = 31249
$ws.Range("E66").Value = "2009"
$ws.Range("F66").Value = 31249
$ws.Range("E67").Value = "2010"
$ws.Range("F67").Value = 31249
$ws.Range("E68").Value = "2011"
$ws.Range("F68").Value = 31249
$ws.Range("E69").Value = "2012"
$ws.Range("F69").Value = 31249
$ws.Range("E70").Value = "2101"
$ws.Range("F70").Value = 31249
$ws.Range("E71").Value = "2102"
$ws.Range("F71").Value = 31249
$ws.Range("E72").Value = "2103"
$ws.Range("F72").Value = 31249
$ws.Range("E73").Value = "2104"
$ws.Range("F73").Value = 31249
$ws.Range("E74").Value = "2105"
$ws.Range("F74").Value = 31249
$ws.Range("E75").Value = "2106"
$ws.Range("F75").Value = 31249
$ws.Range("E76").Value = "2107"
$ws.Range("F76").Value = 31249
$ws.Range("E77").Value = "2108"
$ws.Range("F77").Value = 31249
$ws.Range("E78").Value = "2109"
$ws.Range("F78").Value = 31249
$ws.Range("E79").Value = "2110"
$ws.Range("F79").Value = 31249
$ws.Range("E80").Value = "2111"
$ws.Range("F80").Value = 31249
$ws.Range("E81").Value = "2112"
$ws.Range("F81").Value = 31249
$ws.Range("E82").Value = "2201"
$ws.Range("F82").Value = 31249
$ws.Range("E83").Value = "2202"
$ws.Range("F83").Value = 31249
$ws.Range("E84").Value = "2203"
$ws.Range("F84").Value = 31249
$ws.Range("E85").Value = "2204"
$ws.Range("F85").Value = 31249
$ws.Range("E86").Value = "2205"
$ws.Range("F86").Value = 31249
$ws.Range("E87").Value = "2206"
$ws.Range("F87").Value = 31249
$ws.Range("E88").Value = "2207"
$ws.Range("F88").Value = 31249
$ws.Range("E89").Value = "2208"
$ws.Range("F89").Value = 31249
$ws.Range("E90").Value = "2209"
$ws.Range("F90").Value = 31249
$ws.Range("E91").Value = "2210"
$ws.Range("F91").Value = 31249
$ws.Range("E92").Value = "2211"
$ws.Range("F92").Value = 31249
$ws.Range("E93").Value = "2212"
$ws.Range("F93").Value = 31249
$ws.Range("E94").Value = "2301"
$ws.Range("F94").Value = 31249
$ws.Range("E95").Value = "2302"
$ws.Range("F95").Value = 31249
$ws.Range("E96").Value = "2303"
$ws.Range("F96").Value = 31249
$ws.Range("E97").Value = "2304"
$ws.Range("F97").Value = 31249
$ws.Range("E98").Value = "2305"
$ws.Range("F98").Value = 31249
$ws.Range("E99").Value = "2306"
$ws.Range("F99").Value = 31249
$ws.Range("E100").Value = "2307"
$ws.Range("F100").Value = 31249
$ws.Range("E101").Value = "2308"
$ws.Range("F101").Value = 31249
$ws.Range("E102").Value = "2309"
$ws.Range("F102").Value = 31249
$ws.Range("E103").Value = "2310"
$ws.Range("F103").Value = 31249
$ws.Range("E104").Value = "2311"
$ws.Range("F104").Value = 31249
$ws.Range("E105").Value = "2312"
$ws.Range("F105").Value = 31249
$ws.Range("E106").Value = "2401"
$ws.Range("F106").Value = 31249
$ws.Range("E107").Value = "2402"
$ws.Range("F107").Value = 31249
$ws.Range("E108").Value = "2403"
$ws.Range("F108").Value = 31249
$ws.Range("E109").Value = "2404"
$ws.Range("F109").Value = 31249
$ws.Range("E110").Value = "2405"
$ws.Range("F110").Value = 31249
$ws.Range("E111").Value = "2406"
$ws.Range("F111").Value = 31249
$ws.Range("E112").Value = "2407"
$ws.Range("F112").Value = 31249
$ws.Range("E113").Value = "2408"
$ws.Range("F113").Value = 31249
$ws.Range("E114").Value = "2409"
$ws.Range("F114").Value = 31249
$ws.Range("E115").Value = "2410"
$ws.Range("F115").Value = 31249
$ws.Range("E116").Value = "2411"
$ws.Range("F116").Value = 31249
$ws.Range("E117").Value = "2412"
$ws.Range("F117").Value = 31249
$ws.Range("E118").Value = "2501"
$ws.Range("F118").Value = 31249
$ws.Range("E119").Value = "2502"
$ws.Range("F119").Value = 31249
$ws.Range("E120").Value = "2503"
$ws.Range("F120").Value = 31249
$ws.Range("E121").Value = "2504"
$ws.Range("F121").Value = 31249
$ws.Range("E122").Value = "2505"
$ws.Range("F122").Value = 31249
$ws.Range("E123").Value = "2506"
$ws.Range("F123").Value = 31249
$ws.Range("E124").Value = "2507"
$ws.Range("F124").Value = 31249
$ws.Range("B125").Value = "CC"
$ws.Range("C125").Value = "1143345156"
$ws.Range("D125").Value = "LUIS DONALDO VELAIDES LONDOÑO"
$ws.Range("G125").Value = 781242
$ws.Range("E125").Value = "2508"
$ws.Range("F125").Value = 31249
